# Update automatico via Actualizar 02-09-2021 14-26-34
# Appends new IGPA daily values (rows 760-768) to the IGPA sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IGPA")

# New data rows: date (serial), value (number or "--" for missing data)
$newRows = @(
    @{ Row = 760; Date = 44224; Value = 22247.12 },
    @{ Row = 761; Date = 44225; Value = 21681.02 },
    @{ Row = 762; Date = 44226; Value = "--" },
    @{ Row = 763; Date = 44227; Value = "--" },
    @{ Row = 764; Date = 44228; Value = 22113.41 },
    @{ Row = 765; Date = 44229; Value = 22313.52 },
    @{ Row = 766; Date = 44230; Value = 22245.16 },
    @{ Row = 767; Date = 44231; Value = 22305.07 },
    @{ Row = 768; Date = 44232; Value = 22478.25 }
)

# Template rows already present in the sheet to copy cell formatting from:
#  - row 759: ordinary "date / numeric value" row
#  - row 755: "date / missing value (--)" row
$numericTemplate = $ws.Range("A759:B759")
$missingTemplate = $ws.Range("A755:B755")

foreach ($item in $newRows) {
    $r = $item.Row
    $destRange = $ws.Range("A" + $r + ":B" + $r)

    if ($item.Value -eq "--") {
        $missingTemplate.Copy()
    } else {
        $numericTemplate.Copy()
    }
    $destRange.PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value = $item.Date
    $ws.Cells.Item($r, 2).Value = $item.Value
}

$excel.CutCopyMode = $false

# Extend the "IGPA" defined name to cover the newly added rows.
$wb.Names.Item("IGPA").RefersTo = "=IGPA!`$A`$1:`$B`$768"

# Move the active selection to the last populated cell, matching the
# workbook's saved view state.
$ws.Range("B768").Select() | Out-Null
